# Cluster analysis Fe-number Frost
# - didn't include the SiO2, only the 2 variables
# - clear difference between area1 and rest of areas
#
# The sheet originally held 3 groups (rows 2-4): (0,4), (1,2), (2,1).
# Re-running the clustering with only 2 variables collapses this to
# 2 groups (rows 2-3): (1,4), (0,3) -- and removes the now-unused 3rd row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving group counts
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 4
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 3

# Remove the row that no longer corresponds to a cluster
$ws.Rows.Item(4).Delete()
